$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "¿Terminada?"
$ws.Range("G4").Select()
